$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A62").Value = "Gabriele Gottardi"
$ws.Range("B62").Value = "Federico Zoller | GREP"
$ws.Range("C62").Value = "Andrea  Roveda  | Pinguini Trentini"
$ws.Range("D62").Value = "Luca Perenzoni | CGB Gamberoni"
$ws.Range("E62").Value = "Michele Merighi | Clitoriders"
$ws.Range("F62").Value = "Alessio  Giordano  | FC Schalke 104"
